$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("EkranDatabase")
$ws2 = $wb.Worksheets.Item("ServerDBKurulum")

# Insert a new blank row at position 22 on sheet2 - this shifts the existing
# "Klepe No ve Çıkış No Durum" row (22) down to row 23, preserving its content.
$ws2.Rows.Item(22).Insert()

# --- Write new text values in the same order they were originally authored, ---
# --- so that shared-string table indices line up with the source workbook.  ---

# 1) New "Klepe Harita Durum" row on sheet2 (row 22)
$ws2.Range("C22").Value = "Klepe Harita Durum(ok veya null)"

# 2) New "Ped Harita Onay" row on sheet1 (row 19)
$ws1.Range("C19").Value = "Ped Harita Onay"

# 3) New "Ped ve Çıkış No Onay" row on sheet1 (row 20)
$ws1.Range("C20").Value = "Ped ve Çıkış No Onay"

# 4) Update existing row 18 on sheet1 with the Açma/Kapama çıkış columns
$ws1.Range("E18").Value = "cikisNoAc'lar"
$ws1.Range("F18").Value = "cikisNoKapa'lar"

# 5) Update the (now shifted) Klepe row 23 on sheet2 with Açma/Kapama çıkış columns
$ws2.Range("E23").Value = "çıkışNoAc'lar"
$ws2.Range("F23").Value = "çıkışNoKapa'lar"

# 6) New "Ped Harita Durum" row on sheet2 (row 24)
$ws2.Range("C24").Value = "Ped Harita Durum(ok veya null)"

# 7) New "Ped No ve Çıkış No Durum" row on sheet2 (row 25)
$ws2.Range("C25").Value = "Ped No ve Çıkış No Durum(ok veya null)"

# 8) "Klepe Harita Konumları" label for row 22 on sheet2
$ws2.Range("D22").Value = "Klepe Harita Konumları"

# 9) "Ped Harita Konumları" label for row 24 on sheet2
$ws2.Range("D24").Value = "Ped Harita Konumları"

# --- Fill in the remaining numeric / reused-text cells for each new/changed row ---

# Sheet1 row 19
$ws1.Range("A19").Value = 17
$ws1.Range("B19").Value = 18
$ws1.Range("D19").Value = "klepeHaritalar"
$ws1.Range("E19").Value = 0
$ws1.Range("F19").Value = 0

# Sheet1 row 20
$ws1.Range("A20").Value = 18
$ws1.Range("B20").Value = 19
$ws1.Range("D20").Value = "klepeNo'lar"
$ws1.Range("E20").Value = "cikisNo'lar"
$ws1.Range("F20").Value = 0

# Sheet2 row 22
$ws2.Range("A22").Value = 20
$ws2.Range("B22").Value = 21
$ws2.Range("E22").Value = 0
$ws2.Range("F22").Value = 0

# Sheet2 row 23 index renumber (content/text already shifted down by the insert)
$ws2.Range("A23").Value = 21
$ws2.Range("B23").Value = 22

# Sheet2 row 24
$ws2.Range("A24").Value = 22
$ws2.Range("B24").Value = 23
$ws2.Range("E24").Value = 0
$ws2.Range("F24").Value = 0

# Sheet2 row 25
$ws2.Range("A25").Value = 23
$ws2.Range("B25").Value = 24
$ws2.Range("D25").Value = "fanNo'lar"
$ws2.Range("E25").Value = "çıkışNo'lar"
$ws2.Range("F25").Value = 0

# Restore the selections shown in the saved workbook
$ws1.Range("F19").Select()
$ws2.Range("C31").Select()
